$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("effort")

# Update existing row 16: overrun effort changes from 3 to 5
$ws.Range("B16").Value = 5

# Add new row 17 and 18: copy format of the A-column date cell first
$ws.Range("A16").Copy()
$ws.Range("A17:A18").PasteSpecial(-4122)

$ws.Range("A17").Value = 41185
$ws.Range("B17").Value = 2.25
$ws.Range("D17").Value = "Implementation task overrun and stack usage, not tested yet"

$ws.Range("A18").Value = 41186
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = "Implementation application interrupts, not tested yet. Testing of setEvent/waitForEvent"

# Update selection to match diff (active cell C18)
$ws.Range("C18").Select()
